$d = $word.ActiveDocument

$d.Content.Find.Execute("66×37=", $true, $false, $false, $false, $false, $true, 1, $false, "80×92=", 2) | Out-Null
$d.Content.Find.Execute("20×99=", $true, $false, $false, $false, $false, $true, 1, $false, "96×37=", 2) | Out-Null
$d.Content.Find.Execute("80×42=", $true, $false, $false, $false, $false, $true, 1, $false, "46×21=", 2) | Out-Null
$d.Content.Find.Execute("18×49=", $true, $false, $false, $false, $false, $true, 1, $false, "25×23=", 2) | Out-Null
$d.Content.Find.Execute("56×39=", $true, $false, $false, $false, $false, $true, 1, $false, "73×75=", 2) | Out-Null
$d.Content.Find.Execute("98×72=", $true, $false, $false, $false, $false, $true, 1, $false, "92×60=", 2) | Out-Null
$d.Content.Find.Execute("30×75=", $true, $false, $false, $false, $false, $true, 1, $false, "28×33=", 2) | Out-Null
$d.Content.Find.Execute("49×89=", $true, $false, $false, $false, $false, $true, 1, $false, "55×19=", 2) | Out-Null
$d.Content.Find.Execute("83×59=", $true, $false, $false, $false, $false, $true, 1, $false, "83×23=", 2) | Out-Null
$d.Content.Find.Execute("74×54=", $true, $false, $false, $false, $false, $true, 1, $false, "74×21=", 2) | Out-Null
$d.Content.Find.Execute("35×92=", $true, $false, $false, $false, $false, $true, 1, $false, "79×24=", 2) | Out-Null
$d.Content.Find.Execute("80×44=", $true, $false, $false, $false, $false, $true, 1, $false, "95×87=", 2) | Out-Null
$d.Content.Find.Execute("96×58=", $true, $false, $false, $false, $false, $true, 1, $false, "19×28=", 2) | Out-Null
$d.Content.Find.Execute("48×59=", $true, $false, $false, $false, $false, $true, 1, $false, "88×82=", 2) | Out-Null
$d.Content.Find.Execute("55×57=", $true, $false, $false, $false, $false, $true, 1, $false, "46×56=", 2) | Out-Null
$d.Content.Find.Execute("45×69=", $true, $false, $false, $false, $false, $true, 1, $false, "46×97=", 2) | Out-Null
$d.Content.Find.Execute("38×58=", $true, $false, $false, $false, $false, $true, 1, $false, "50×86=", 2) | Out-Null
$d.Content.Find.Execute("66×45=", $true, $false, $false, $false, $false, $true, 1, $false, "20×34=", 2) | Out-Null
$d.Content.Find.Execute("15×53=", $true, $false, $false, $false, $false, $true, 1, $false, "93×85=", 2) | Out-Null
$d.Content.Find.Execute("37×86=", $true, $false, $false, $false, $false, $true, 1, $false, "55×32=", 2) | Out-Null
$d.Content.Find.Execute("95×99=", $true, $false, $false, $false, $false, $true, 1, $false, "45×23=", 2) | Out-Null
$d.Content.Find.Execute("61×62=", $true, $false, $false, $false, $false, $true, 1, $false, "84×71=", 2) | Out-Null
$d.Content.Find.Execute("97×22=", $true, $false, $false, $false, $false, $true, 1, $false, "59×88=", 2) | Out-Null
$d.Content.Find.Execute("11×85=", $true, $false, $false, $false, $false, $true, 1, $false, "86×50=", 2) | Out-Null
$d.Content.Find.Execute("15×76=", $true, $false, $false, $false, $false, $true, 1, $false, "59×28=", 2) | Out-Null
